# Update the NIEM 3.1 mapping paths for the Driver License identification
# elements: "DriverLicenseCardIdentification" -> "DriverLicenseIdentification".
# (Commit: Added "DriverLicenseCardIdentification" to NIEM Subset)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C18").Value = "cs-req-doc:CustodySearchRequest/nc:Person/j:PersonAugmentation/j:DriverLicense/j:DriverLicenseIdentification/nc:IdentificationID"
$ws.Range("C19").Value = "cs-req-doc:CustodySearchRequest/nc:Person/j:PersonAugmentation/j:DriverLicense/j:DriverLicenseIdentification/nc:IdentificationSourceText"

# Mirror the author's final cursor position / selection in the saved view.
$ws.Range("C16").Select()
